$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '46.194.18'
$ws.Range('E2').Value = '  -0.01%  '
Set-TextValue 'D3' '2.600.06'
$ws.Range('E3').Value = '  +6.64%  '
$ws.Range('E4').Value = '  +0.14%  '
Set-TextValue 'D5' '307.55'
$ws.Range('E5').Value = '  +3.75%  '
Set-TextValue 'D6' '99.76'
$ws.Range('E6').Value = '  +3.89%  '
Set-TextValue 'D7' '0.602'
$ws.Range('E7').Value = '  +5.63%  '
$ws.Range('E8').Value = '  +0.16%  '
Set-TextValue 'D9' '0.580'
$ws.Range('E9').Value = '  +14.05%  '
Set-TextValue 'D10' '39.21'
$ws.Range('E10').Value = '  +11.09%  '
Set-TextValue 'D11' '54.44'
$ws.Range('E11').Value = '  +1.22%  '
Set-TextValue 'D12' '0.0844'
$ws.Range('E12').Value = '  +7.52%  '
Set-TextValue 'D13' '8.16'
$ws.Range('E13').Value = '  +14.51%  '
Set-TextValue 'D14' '2.986.67'
$ws.Range('E14').Value = '  +6.71%  '
Set-TextValue 'D15' '0.106'
$ws.Range('E15').Value = '  +1.20%  '
Set-TextValue 'D16' '2.600.98'
$ws.Range('E16').Value = '  +7.82%  '
Set-TextValue 'D17' '0.923'
$ws.Range('E17').Value = '  +9.06%  '
Set-TextValue 'D18' '15.03'
$ws.Range('E18').Value = '  +5.70%  '
Set-TextValue 'D19' '46.287.39'
$ws.Range('E19').Value = '  +0.66%  '
$ws.Range('E20').Value = '  +6.84%  '
Set-TextValue 'D21' '13.00'
$ws.Range('E21').Value = '  +3.78%  '
Set-TextValue 'D22' '6.75'
$ws.Range('E22').Value = '  +8.69%  '
Set-TextValue 'D23' '71.43'
$ws.Range('E23').Value = '  +5.87%  '
Set-TextValue 'D24' '273.78'
$ws.Range('E24').Value = '  +12.20%  '
Set-TextValue 'D25' '3.03'
$ws.Range('E25').Value = '  +8.30%  '
Set-TextValue 'D26' '2.17'
$ws.Range('E26').Value = '  +11.05%  '
Set-TextValue 'D27' '29.65'
$ws.Range('E27').Value = '  +38.83%  '
Set-TextValue 'D28' '1.00'
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('E29').Value = '  +0.19%  '
Set-TextValue 'D30' '10.60'
$ws.Range('E30').Value = '  +8.53%  '
Set-TextValue 'D31' '2.32'
$ws.Range('E31').Value = '  +4.24%  '
Set-TextValue 'D32' '39.11'
$ws.Range('E32').Value = '  -0.45%  '
Set-TextValue 'D33' '6.33'
$ws.Range('E33').Value = '  +14.38%  '
Set-TextValue 'D34' '3.62'
$ws.Range('E34').Value = '  -6.73%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D35' '0.0842'
$ws.Range('E35').Value = '  +9.16%  '
$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D36' '2.83'
$ws.Range('E36').Value = '  +2.19%  '
Set-TextValue 'D37' '2.20'
$ws.Range('E37').Value = '  +8.74%  '
Set-TextValue 'D38' '150.30'
$ws.Range('E38').Value = '  +1.18%  '
$ws.Range('E39').Value = '  +6.08%  '
$ws.Range('E40').Value = '  +5.38%  '
Set-TextValue 'D41' '23.07'
$ws.Range('E41').Value = '  +40.18%  '
Set-TextValue 'D42' '16.03'
$ws.Range('E42').Value = '  +7.77%  '
$ws.Range('B43').Value = 'NEARProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D43' '3.63'
$ws.Range('E43').Value = '  +11.14%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D44' '0.0330'
$ws.Range('E44').Value = '  +9.52%  '
Set-TextValue 'D45' '4.08'
$ws.Range('E45').Value = '  +6.68%  '
Set-TextValue 'D46' '2.160.77'
$ws.Range('E46').Value = '  +8.61%  '
$ws.Range('E47').Value = '  -0.08%  '
Set-TextValue 'D48' '93.73'
$ws.Range('E48').Value = '  +4.14%  '
Set-TextValue 'D49' '9.55'
$ws.Range('E49').Value = '  +10.87%  '
Set-TextValue 'D50' '109.61'
$ws.Range('E50').Value = '  +7.95%  '
Set-TextValue 'D51' '1.78'
$ws.Range('E51').Value = '  -2.26%  '
